$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.183.81"
$ws.Range("E2").Value = "  -2.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.819.31"
$ws.Range("E3").Value = "  -2.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -1.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.00"
$ws.Range("E5").Value = "  -2.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("E6").Value = "  -1.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4252"
$ws.Range("E7").Value = "  -2.69%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3675"
$ws.Range("E8").Value = "  -2.91%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07223"
$ws.Range("E9").Value = "  -2.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8600"
$ws.Range("E10").Value = "  -2.94%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.94"
$ws.Range("E11").Value = "  -3.48%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.818.55"
$ws.Range("E12").Value = "  -2.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.667"
$ws.Range("E13").Value = "  -1.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.07101"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.302"
$ws.Range("E15").Value = "  -3.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.00"
$ws.Range("E16").Value = "  -0.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008856"
$ws.Range("E18").Value = "  -2.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("E20").Value = "  -2.94%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.227.30"
$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.133"
$ws.Range("E22").Value = "  -2.65%  "

$ws.Range("E23").Value = "  -3.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.046.19"
$ws.Range("E24").Value = "  -1.63%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.004"
$ws.Range("E25").Value = "  -1.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.48"
$ws.Range("E26").Value = "  -2.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.27"
$ws.Range("E27").Value = "  -2.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.112"
$ws.Range("E28").Value = "  +5.39%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.215"
$ws.Range("E29").Value = "  -4.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.07"
$ws.Range("E30").Value = "  -3.79%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08867"
$ws.Range("E31").Value = "  -1.89%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.192"
$ws.Range("E32").Value = "  -3.65%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7557"
$ws.Range("E33").Value = "  -2.36%  "

$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.807"
$ws.Range("E35").Value = "  -6.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").Value = "  -1.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.111"
$ws.Range("E37").Value = "  -2.70%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01970"
$ws.Range("E38").Value = "  -0.58%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05263"
$ws.Range("E39").Value = "  -1.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.121"
$ws.Range("E40").Value = "  +1.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.862"
$ws.Range("E41").Value = "  -0.57%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1687"
$ws.Range("E42").Value = "  +0.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5031"
$ws.Range("E43").Value = "  -3.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.597"
$ws.Range("E44").Value = "  -1.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.63"
$ws.Range("E45").Value = "  -1.13%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "106.74"
$ws.Range("E46").Value = "  -3.40%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4732"
$ws.Range("E47").Value = "  -0.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  -1.30%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06367"
$ws.Range("E49").Value = "  -1.75%  "

$ws.Range("E50").Value = "  -3.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.801"
$ws.Range("E51").Value = "  -2.90%  "
